$p = $ppt.ActivePresentation

# Slide 8 (SlideID 690, hidden) "The Basics of Input-Output Analysis" is the
# slide that was hidden (per the discussion thread in its comments) and is
# now being permanently removed from the deck.
$s = $p.Slides.Item(8)

# Remove the comment thread (top-level comment + its replies) attached to
# this slide first, so the comments part is fully cleaned up rather than
# left as an orphan when the slide itself is deleted.
for ($i = $s.Comments.Count; $i -ge 1; $i--) {
    $c = $s.Comments.Item($i)
    for ($j = $c.Replies.Count; $j -ge 1; $j--) {
        $c.Replies.Item($j).Delete()
    }
    $c.Delete()
}

# Deleting the slide also removes its associated notes slide, and shifts
# the following slide (SlideID 682, "Suggested Readings:") up into its
# place.
$s.Delete()
